$wb = $excel.ActiveWorkbook

# Update zh-cn sheet (row 4: Correspond Handoff Datetime / Correspond Handback DateTime)
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E4").Value = "2016-03-14 06:33:21"
$wsZh.Range("H4").Value = "2016-03-14 06:33:38"

# Update de-de sheet (row 4: Correspond Handoff Datetime / Correspond Handback DateTime)
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E4").Value = "2016-03-14 06:33:24"
$wsDe.Range("H4").Value = "2016-03-14 06:33:44"
